$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 322.13333
$ws.Range("I33").Value = 232.44444
$ws.Range("K33").Value = 232.44444
$ws.Range("M33").Value = -3.444439999999986
$ws.Range("H43").Value = 4250
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 7500
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 7500
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -7638
$ws.Range("H64").Value = 31368.25
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 31368.25
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H107").Value = 1500.619
$ws.Range("I107").Value = 1181.2
$ws.Range("K107").Value = 1181.2
$ws.Range("M107").Value = 738.8
$ws.Range("H109").Value = 105000
$ws.Range("J109").Value = 105000
$ws.Range("L109").Value = 105000
$ws.Range("N109").Value = -107774
$ws.Range("H112").Value = 6087.8237
$ws.Range("I112").Value = 2399.5
$ws.Range("K112").Value = 7198.5
$ws.Range("M112").Value = -6090.5
$ws.Range("H138").Value = 18863.639
$ws.Range("J138").Value = 32466.354
$ws.Range("L138").Value = 97399.06200000001
$ws.Range("N138").Value = -107679.062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20971.527
$ws.Range("I32").Value = 22814.459
$ws.Range("K32").Value = 22814.459
$ws.Range("M32").Value = -22527.459
$ws.Range("H45").Value = 2995.75
$ws.Range("I45").Value = 2276.3333
$ws.Range("K45").Value = 2276.3333
$ws.Range("M45").Value = -1899.3333
$ws.Range("H61").Value = 5812.615
$ws.Range("J61").Value = 15451.375
$ws.Range("L61").Value = 15451.375
$ws.Range("N61").Value = -15875.375
$ws.Range("H74").Value = 306609.66
$ws.Range("I74").Value = 600879.2
$ws.Range("K74").Value = 600879.2
$ws.Range("M74").Value = -600005.2
$ws.Range("H77").Value = 306609.66
$ws.Range("I77").Value = 600879.2
$ws.Range("K77").Value = 3004396
$ws.Range("M77").Value = -3000028
$ws.Range("H110").Value = 27968.455
$ws.Range("I110").Value = 30230.133
$ws.Range("K110").Value = 30230.133
$ws.Range("M110").Value = -28185.133
$ws.Range("H122").Value = 1679.1034
$ws.Range("I122").Value = 1463.3334
$ws.Range("J122").Value = 2714.8
$ws.Range("K122").Value = 4390.0002
$ws.Range("L122").Value = 8144.400000000001
$ws.Range("M122").Value = -1940.0002
$ws.Range("N122").Value = -13044.4
$ws.Range("H132").Value = 1868.88
$ws.Range("I132").Value = 1320.1428
$ws.Range("K132").Value = 3960.4284
$ws.Range("M132").Value = -1430.4284
$ws.Range("H135").Value = 91525
$ws.Range("J135").Value = 91525
$ws.Range("L135").Value = 91525
$ws.Range("N135").Value = -101665
$ws.Range("H136").Value = 5812.615
$ws.Range("J136").Value = 15451.375
$ws.Range("L136").Value = 46354.125
$ws.Range("N136").Value = -51454.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6993.7646
$ws.Range("I20").Value = 8409.519
$ws.Range("J20").Value = 1533
$ws.Range("K20").Value = 8409.519
$ws.Range("L20").Value = 1533
$ws.Range("M20").Value = -8162.519
$ws.Range("N20").Value = -2027
$ws.Range("H86").Value = 2284.4167
$ws.Range("I86").Value = 2056.1428
$ws.Range("J86").Value = 2604
$ws.Range("K86").Value = 2056.1428
$ws.Range("L86").Value = 2604
$ws.Range("M86").Value = -933.1428000000001
$ws.Range("N86").Value = -4850
$ws.Range("H89").Value = 2284.4167
$ws.Range("I89").Value = 2056.1428
$ws.Range("J89").Value = 2604
$ws.Range("K89").Value = 10280.714
$ws.Range("L89").Value = 13020
$ws.Range("M89").Value = -4664.714
$ws.Range("N89").Value = -24252
$ws.Range("H134").Value = 3353.0588
$ws.Range("I134").Value = 3083.8572
$ws.Range("K134").Value = 9251.571599999999
$ws.Range("M134").Value = -6716.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50013950
$ws.Range("I31").Value = 50013950
$ws.Range("K31").Value = 50013950
$ws.Range("M31").Value = -50013655
$ws.Range("H34").Value = 50013950
$ws.Range("I34").Value = 50013950
$ws.Range("K34").Value = 50013950
$ws.Range("M34").Value = -50013748
$ws.Range("H58").Value = 1503.0322
$ws.Range("I58").Value = 1475.7142
$ws.Range("K58").Value = 1475.7142
$ws.Range("M58").Value = -1272.7142
$ws.Range("H105").Value = 2351
$ws.Range("I105").Value = 1688.75
$ws.Range("K105").Value = 1688.75
$ws.Range("M105").Value = 58.25
$ws.Range("H132").Value = 84482.414
$ws.Range("I132").Value = 84482.414
$ws.Range("K132").Value = 253447.242
$ws.Range("M132").Value = -250917.242
$ws.Range("H134").Value = 2192.6296
$ws.Range("I134").Value = 1535.2858
$ws.Range("K134").Value = 4605.857400000001
$ws.Range("M134").Value = -2070.857400000001
$ws.Range("H136").Value = 1503.0322
$ws.Range("I136").Value = 1475.7142
$ws.Range("K136").Value = 4427.142599999999
$ws.Range("M136").Value = -1877.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 3147.5
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 3530
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 10590
$ws.Range("M48").Value = -5750
$ws.Range("N48").Value = -11090
$ws.Range("H56").Value = 6830.6665
$ws.Range("I56").Value = 6830.6665
$ws.Range("K56").Value = 6830.6665
$ws.Range("M56").Value = -6300.6665
$ws.Range("H131").Value = 429090.7
$ws.Range("I131").Value = 429090.7
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1287272.1
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1282232.1
$ws.Range("N131").ClearContents()
$ws.Range("H137").Value = 4619.8
$ws.Range("I137").Value = 4366.3335
$ws.Range("K137").Value = 13099.0005
$ws.Range("M137").Value = -7999.000499999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8778
$ws.Range("J41").Value = 18950
$ws.Range("L41").Value = 18950
$ws.Range("N41").Value = -19660
$ws.Range("H102").Value = 1949.5
$ws.Range("I102").Value = 1326.7222
$ws.Range("J102").Value = 3817.8333
$ws.Range("K102").Value = 1326.7222
$ws.Range("L102").Value = 3817.8333
$ws.Range("M102").Value = 295.2778000000001
$ws.Range("N102").Value = -7061.8333
$ws.Range("H126").Value = 3374.75
$ws.Range("I126").Value = 2499.6667
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 7499.000100000001
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -5029.000100000001
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2863.923
$ws.Range("I122").Value = 2854.0417
$ws.Range("J122").Value = 2982.5
$ws.Range("K122").Value = 8562.125100000001
$ws.Range("L122").Value = 8947.5
$ws.Range("M122").Value = -6112.125100000001
$ws.Range("N122").Value = -13847.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 11179.8
$ws.Range("I3").Value = 7949.5
$ws.Range("K3").Value = 7949.5
$ws.Range("M3").Value = -7835.5
$ws.Range("H45").Value = 16105.454
$ws.Range("J45").Value = 16105.454
$ws.Range("L45").Value = 16105.454
$ws.Range("N45").Value = -17087.454
$ws.Range("H107").Value = 875.4
$ws.Range("I107").Value = 831.6
$ws.Range("J107").Value = 963
$ws.Range("K107").Value = 2494.8
$ws.Range("L107").Value = 2889
$ws.Range("M107").Value = -574.8000000000002
$ws.Range("N107").Value = -6729
$ws.Range("H132").Value = 32385.479
$ws.Range("I132").Value = 36693.6
$ws.Range("K132").Value = 110080.8
$ws.Range("M132").Value = -107550.8
